# Add a new "status_label" column right after the existing "statut" column.
# This shifts every existing column from B..I to C..J, which Excel's
# Columns.Insert() does natively (formats/content alike).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("B:B").Insert()
$ws.Range("B1").Value = "status_label"

# Populate the new status_label column for every data row, derived from the
# emoji already stored in column A (🟥 -> rouge, 🟩 -> vert, 🟧 -> orange).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $status = $ws.Cells.Item($r, 1).Value2
    if ($status -eq "🟥") {
        $label = "rouge"
    } elseif ($status -eq "🟩") {
        $label = "vert"
    } elseif ($status -eq "🟧") {
        $label = "orange"
    } else {
        $label = ""
    }
    $ws.Cells.Item($r, 2).Value = $label
}

# Two pairs of rows were also re-ordered by NCTId. The completion_year
# (column E) and eudraCT (column D, always blank) are identical within each
# swapped pair, so only the status/NCTId/title/acronym/results columns need
# to move.

# Swap rows 6 and 7 (NCT02000674 <-> NCT01425866); both keep status "rouge"
# and both already carry the same results_1y/3y/results booleans, so only
# the NCTId/title/acronym need to trade places.
$ws.Range("C6").Value = "NCT01425866"
$ws.Range("F6").Value = "Multicenter Randomized Trial of Structured Educational Intervention at the Community Level in Insufficiently Controlled Patients With Type 2 Diabetes in Reunion Island"
$ws.Range("G6").Value = "ERMIES"

$ws.Range("C7").Value = "NCT02000674"
$ws.Range("F7").Value = "Succinylcholine vs Rocuronium for Prehospital Emergency Intubation : a Randomized Trial"
$ws.Range("G7").Value = "CURASMUR"

# Swap rows 13 and 14 (NCT04459221 [vert] <-> NCT05098925 [rouge]); the
# status marker (and its label) travel with the row this time.
$ws.Range("A13").Value = "🟥"
$ws.Range("B13").Value = "rouge"
$ws.Range("C13").Value = "NCT05098925"
$ws.Range("F13").Value = "Study of Thermoregulatory Processes in Ultra-endurance Runners in a Hot and Humid Environment"
$ws.Range("G13").Value = "ERUPTION-2"
$ws.Range("H13").Value = $false
$ws.Range("I13").Value = $false
$ws.Range("J13").Value = $false

$ws.Range("A14").Value = "🟩"
$ws.Range("B14").Value = "vert"
$ws.Range("C14").Value = "NCT04459221"
$ws.Range("F14").Value = "Study of the Impact of a School Program Combining - Promotion of HPV Vaccination and HPV Vaccine Offer in Middle School - on Adherence to HPV Vaccination in Middle School Students"
$ws.Range("G14").Value = "PROM SSCOL"
$ws.Range("H14").Value = $true
$ws.Range("I14").Value = $true
$ws.Range("J14").Value = $true
